$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1802.2778
$ws.Range("J17").Value = 1802.2778
$ws.Range("L17").Value = 5406.8334
$ws.Range("N17").Value = -5742.8334

# Row 19
$ws.Range("H19").Value = 1264.15
$ws.Range("I19").Value = 1259.0588
$ws.Range("J19").Value = 1293
$ws.Range("K19").Value = 1259.0588
$ws.Range("L19").Value = 1293
$ws.Range("M19").Value = -1084.0588
$ws.Range("N19").Value = -1643

# Row 33
$ws.Range("H33").Value = 248.8
$ws.Range("I33").Value = 313.2857
$ws.Range("J33").Value = 98.333336
$ws.Range("K33").Value = 313.2857
$ws.Range("L33").Value = 98.333336
$ws.Range("M33").Value = -84.28570000000002
$ws.Range("N33").Value = -556.333336

# Row 62
$ws.Range("H62").Value = 2649
$ws.Range("I62").Value = 1849
$ws.Range("K62").Value = 1849
$ws.Range("M62").Value = -1225

# Row 65
$ws.Range("H65").Value = 2649
$ws.Range("I65").Value = 1849
$ws.Range("K65").Value = 9245
$ws.Range("M65").Value = -6125

# Row 107
$ws.Range("H107").Value = 1114.75
$ws.Range("I107").Value = 1086.6364
$ws.Range("J107").Value = 1176.6
$ws.Range("K107").Value = 1086.6364
$ws.Range("L107").Value = 1176.6
$ws.Range("M107").Value = 833.3635999999999
$ws.Range("N107").Value = -5016.6

# Row 116
$ws.Range("H116").Value = 2082
$ws.Range("I116").Value = 2150.8572
$ws.Range("K116").Value = 2150.8572
$ws.Range("M116").Value = 1291.1428

# Row 118
$ws.Range("H118").Value = 847.1
$ws.Range("I118").Value = 847.1
$ws.Range("K118").Value = 2541.3
$ws.Range("M118").Value = -884.3000000000002

# Row 127
$ws.Range("H127").Value = 1783.5
$ws.Range("I127").Value = 1783.5
$ws.Range("K127").Value = 5350.5
$ws.Range("M127").Value = -390.5

$ws = $wb.Worksheets.Item("ARM")
# Row 110
$ws.Range("H110").Value = 7580.2
$ws.Range("I110").Value = 8045.7144
$ws.Range("K110").Value = 8045.7144
$ws.Range("M110").Value = -6000.7144

# Row 113
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

# Row 120
$ws.Range("H120").Value = 33713.5
$ws.Range("J120").Value = 36380
$ws.Range("L120").Value = 36380
$ws.Range("N120").Value = -46056

# Row 122
$ws.Range("H122").Value = 1288.3
$ws.Range("I122").Value = 1288.3
$ws.Range("K122").Value = 3864.9
$ws.Range("M122").Value = -1414.9

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 5662.3335
$ws.Range("I86").Value = 3194.9
$ws.Range("J86").Value = 10597.2
$ws.Range("K86").Value = 3194.9
$ws.Range("L86").Value = 10597.2
$ws.Range("M86").Value = -2071.9
$ws.Range("N86").Value = -12843.2

# Row 89
$ws.Range("H89").Value = 5662.3335
$ws.Range("I89").Value = 3194.9
$ws.Range("J89").Value = 10597.2
$ws.Range("K89").Value = 15974.5
$ws.Range("L89").Value = 52986
$ws.Range("M89").Value = -10358.5
$ws.Range("N89").Value = -64218

# Row 134
$ws.Range("H134").Value = 17999.5
$ws.Range("I134").Value = 17999.5
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 53998.5
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -51463.5
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1687.2916
$ws.Range("I31").Value = 1677.1364
$ws.Range("K31").Value = 1677.1364
$ws.Range("M31").Value = -1382.1364

# Row 34
$ws.Range("H34").Value = 1687.2916
$ws.Range("I34").Value = 1677.1364
$ws.Range("K34").Value = 1677.1364
$ws.Range("M34").Value = -1475.1364

# Row 99
$ws.Range("H99").Value = 3848.2856
$ws.Range("I99").Value = 3320.1667
$ws.Range("K99").Value = 3320.1667
$ws.Range("M99").Value = -1822.1667

# Row 103
$ws.Range("H103").Value = 25501
$ws.Range("I103").Value = 25501
$ws.Range("K103").Value = 25501
$ws.Range("M103").Value = -24329

# Row 122
$ws.Range("H122").Value = 3651.3667
$ws.Range("I122").Value = 3674.6
$ws.Range("K122").Value = 11023.8
$ws.Range("M122").Value = -8573.799999999999

# Row 126
$ws.Range("H126").Value = 3848.2856
$ws.Range("I126").Value = 3320.1667
$ws.Range("K126").Value = 9960.500100000001
$ws.Range("M126").Value = -7490.500100000001

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 836.9048
$ws.Range("I5").Value = 814.4286
$ws.Range("K5").Value = 2443.2858
$ws.Range("M5").Value = -2331.2858

# Row 7
$ws.Range("H7").Value = 87591.836
$ws.Range("I7").Value = 146039.14
$ws.Range("K7").Value = 438117.42
$ws.Range("M7").Value = -438005.42

# Row 92
$ws.Range("H92").Value = 156.44444
$ws.Range("I92").Value = 194.75
$ws.Range("J92").Value = 125.8
$ws.Range("K92").Value = 584.25
$ws.Range("L92").Value = 377.4
$ws.Range("M92").Value = 663.75
$ws.Range("N92").Value = -2873.4

# Row 113
$ws.Range("H113").Value = 1142.7576
$ws.Range("I113").Value = 652.15
$ws.Range("J113").Value = 1897.5385
$ws.Range("K113").Value = 1956.45
$ws.Range("L113").Value = 5692.6155
$ws.Range("M113").Value = 213.5500000000002
$ws.Range("N113").Value = -10032.6155

# Row 135
$ws.Range("H135").Value = 836.9048
$ws.Range("I135").Value = 814.4286
$ws.Range("K135").Value = 7329.8574
$ws.Range("M135").Value = -4794.8574

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 3635
$ws.Range("I122").Value = 3215.5
$ws.Range("K122").Value = 9646.5
$ws.Range("M122").Value = -7196.5

# Row 126
$ws.Range("H126").Value = 4776.4
$ws.Range("I126").Value = 4443
$ws.Range("K126").Value = 13329
$ws.Range("M126").Value = -10859

# Row 132
$ws.Range("H132").Value = 2434.3333
$ws.Range("I132").Value = 2434.3333
$ws.Range("K132").Value = 7302.999899999999
$ws.Range("M132").Value = -4772.999899999999

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 5120
$ws.Range("J7").Value = 5550
$ws.Range("L7").Value = 5550
$ws.Range("N7").Value = -5774

# Row 18
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

# Row 126
$ws.Range("H126").Value = 5120
$ws.Range("J126").Value = 5550
$ws.Range("L126").Value = 16650
$ws.Range("N126").Value = -21590

# Row 132
$ws.Range("H132").Value = 2336.6667
$ws.Range("I132").Value = 2323.318
$ws.Range("K132").Value = 6969.954000000001
$ws.Range("M132").Value = -4439.954000000001

$ws = $wb.Worksheets.Item("WVR")
# Row 45
$ws.Range("H45").Value = 42555.184
$ws.Range("I45").Value = 31810.8
$ws.Range("J45").Value = 51508.832
$ws.Range("K45").Value = 31810.8
$ws.Range("L45").Value = 51508.832
$ws.Range("M45").Value = -31319.8
$ws.Range("N45").Value = -52490.832

# Row 62
$ws.Range("H62").Value = 5129
$ws.Range("I62").Value = 4538.2856
$ws.Range("K62").Value = 4538.2856
$ws.Range("M62").Value = -3914.2856

# Row 65
$ws.Range("H65").Value = 5129
$ws.Range("I65").Value = 4538.2856
$ws.Range("K65").Value = 22691.428
$ws.Range("M65").Value = -19571.428

# Row 81
$ws.Range("H81").Value = 3271.1738
$ws.Range("I81").Value = 3366.85
$ws.Range("J81").Value = 2633.3333
$ws.Range("K81").Value = 6733.7
$ws.Range("L81").Value = 5266.6666
$ws.Range("M81").Value = -5672.7
$ws.Range("N81").Value = -7388.6666

# Row 84
$ws.Range("H84").Value = 3271.1738
$ws.Range("I84").Value = 3366.85
$ws.Range("J84").Value = 2633.3333
$ws.Range("K84").Value = 33668.5
$ws.Range("L84").Value = 26333.333
$ws.Range("M84").Value = -28364.5
$ws.Range("N84").Value = -36941.333

# Row 96
$ws.Range("H96").Value = 3446.8
$ws.Range("I96").Value = 2745.3333
$ws.Range("J96").Value = 4499
$ws.Range("K96").Value = 2745.3333
$ws.Range("L96").Value = 4499
$ws.Range("M96").Value = -1372.3333
$ws.Range("N96").Value = -7245

# Row 126
$ws.Range("H126").Value = 1922.7778
$ws.Range("I126").Value = 1286.5714
$ws.Range("J126").Value = 4149.5
$ws.Range("K126").Value = 3859.7142
$ws.Range("L126").Value = 12448.5
$ws.Range("M126").Value = -1389.7142
$ws.Range("N126").Value = -17388.5
